# Edit: "More accessible chart for communications method activity for students"
$d = $word.ActiveDocument

# --- Change 1: split "Depending on the amount of time available, the Communications
# Methods activity..." run and insert a collapsed _GoBack bookmark at the split point
# (this also relocates the document's existing _GoBack bookmark away from the end of
# the document, covering change 4 below).
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Text = "Depending on the amount of time available, the Communications"
$found1 = $r1.Find.Execute()
$splitPos = $r1.Start + 20
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Change 2: insert a new green-highlighted list item ("We created a list view...")
# right before the "This assignment..." bullet, and give that original bullet's own
# paragraph mark the green highlight too, matching the XML diff exactly.
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Text = "This assignment, in a normal class setting"
$found2 = $r2.Find.Execute()
$para2 = $r2.Paragraphs(1)
$prange2 = $para2.Range
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="720"/></w:tabs><w:suppressAutoHyphens/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times;Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:highlight w:val="green"/></w:rPr><w:t>We created a list view of the chart in the “Revised” version of the chart files for visually impaired students.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="13"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="720"/></w:tabs><w:suppressAutoHyphens/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times;Times New Roman" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>This assignment, in a normal class setting, will be due on the first day of Unit 2.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $prange2.InsertXML($xml2)

# --- Change 3: merge the two runs of the "Let's look at what kinds of data..." /
# "is not just the content..." paragraph into a single run, moving
# <w:lastRenderedPageBreak/> to the front of the merged run.
$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Find.Text = "Let's look at what kinds"
$found3 = $r3.Find.Execute()
$r3end = $d.Content
$r3end.Find.ClearFormatting()
$r3end.Find.Text = "telephone call and the number."
$foundEnd = $r3end.Find.Execute()
$full3 = $d.Range($r3.Start, $r3end.End)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times;Times New Roman" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t>Let''s look at what kinds of data you ''give off'' when using the different forms of communication. For each of the following examples, fill in which method you would choose for the given scenario and why (You should already have completed that part.). Keep in mind that “data” here is not just the content you communicate (what you say or write) but could also refer to details like the time of a telephone call and the number.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $full3.InsertXML($xml3)

Write-Output "changes applied"
